# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" sheet (cloned from "2022-Q3" so it keeps the
#   same header/row styling) positioned right after "总计" and before
#   "2022-Q3".
# - Update the "总计" (summary) sheet: new first data row for 2022-Q4,
#   shifting the existing 2022-Q3 / 2022-Q2 rows down by one.

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (t="s"/inlineStr) instead of
    # being auto-coerced to a number - several columns here hold
    # numeric-looking strings (fund codes with leading zeros, percentages,
    # etc.) that must stay text, matching the original workbook.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop back to the default style so we don't leave a stray
    # quotePrefix-flavoured style on the cell.
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step 1: update "总计" - insert the 2022-Q4 row above the existing data
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 4 doesn't exist yet - clone the index-column style from row 3 first.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 = old Row 3 data (2022-Q2), now at index 2
$total.Range("A4").Value = 2
Set-TextValue $total.Range("B4") "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 2.78

# Row 3 = old Row 2 data (2022-Q3), now at index 1
$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2022-Q3"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 2.09

# Row 2 = new 2022-Q4 data, at index 0
$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 1.29

# ---------------------------------------------------------------
# Step 2: create the "2022-Q4" sheet, positioned before "2022-Q3"
# ---------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcIndex = $srcQ3.Index
$srcQ3.Copy($srcQ3)
$q4 = $wb.Worksheets.Item($srcIndex)
$q4.Name = "2022-Q4"

# The source sheet has 6 data rows; 2022-Q4 only needs 4, so drop 2.
$q4.Rows(6).Delete()
$q4.Rows(6).Delete()

# Row 2
Set-TextValue $q4.Range("B2") "013776"
Set-TextValue $q4.Range("C2") "中泰兴为价值精选混合A"
Set-TextValue $q4.Range("D2") "12.01"
Set-TextValue $q4.Range("E2") "92.98"
Set-TextValue $q4.Range("F2") "4.90"
Set-TextValue $q4.Range("G2") "0.5885"
$q4.Range("H2").Value = 8

# Row 3
Set-TextValue $q4.Range("B3") "010728"
Set-TextValue $q4.Range("C3") "中泰兴诚价值一年持有期混合A"
Set-TextValue $q4.Range("D3") "6.60"
Set-TextValue $q4.Range("E3") "92.29"
Set-TextValue $q4.Range("F3") "5.04"
Set-TextValue $q4.Range("G3") "0.3326"
$q4.Range("H3").Value = 10

# Row 4
Set-TextValue $q4.Range("B4") "013777"
Set-TextValue $q4.Range("C4") "中泰兴为价值精选混合C"
Set-TextValue $q4.Range("D4") "6.16"
Set-TextValue $q4.Range("E4") "92.98"
Set-TextValue $q4.Range("F4") "4.90"
Set-TextValue $q4.Range("G4") "0.3018"
$q4.Range("H4").Value = 8

# Row 5
Set-TextValue $q4.Range("B5") "010729"
Set-TextValue $q4.Range("C5") "中泰兴诚价值一年持有期混合C"
Set-TextValue $q4.Range("D5") "1.24"
Set-TextValue $q4.Range("E5") "92.29"
Set-TextValue $q4.Range("F5") "5.04"
Set-TextValue $q4.Range("G5") "0.0625"
$q4.Range("H5").Value = 10

# ---------------------------------------------------------------
# Step 3: copying a sheet makes the new copy the active tab - restore
# "2022-Q2" (the originally-selected sheet) as active.
# ---------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
